# Update the "loading_percent" results for the 380 kV case (column F, rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 16.86991607391245
    3  = 15.89584955866815
    4  = 15.2699775710849
    5  = 15.008197319934
    6  = 14.96433081551589
    7  = 15.26647399323134
    8  = 16.53996406344772
    9  = 19.0027458068253
    10 = 20.67494806633232
    11 = 21.39172256362241
    12 = 21.65686569030329
    13 = 21.60004134736741
    14 = 21.4136618050453
    15 = 21.29868154950794
    16 = 20.62722412089977
    17 = 20.20408069597326
    18 = 19.95656407809801
    19 = 19.87204792380568
    20 = 20.24955283636154
    21 = 21.46857628470571
    22 = 22.22866616901554
    23 = 21.82633154458858
    24 = 20.22900810905285
    25 = 18.34778573295695
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
